$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, shifting existing rows 180-220 down to 181-221.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new price record.
$ws.Range("A180").Value() = 7
$ws.Range("B180").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C180").Value() = "Ñuble"
$ws.Range("D180").Value() = "2021-11-08"
$ws.Range("E180").Value() = 16
$ws.Range("F180").Value() = 100114001
$ws.Range("G180").Value() = "Papa"
$ws.Range("H180").Value() = "Patagonia"
$ws.Range("I180").Value() = "1a (guarda)"
$ws.Range("J180").Value() = 400
$ws.Range("K180").Value() = 7000
$ws.Range("L180").Value() = 8000
$ws.Range("M180").Value() = 7500
$ws.Range("N180").Value() = "$/saco 25 kilos"
$ws.Range("O180").Value() = "Provincia de Diguillín"
$ws.Range("P180").Value() = 300
$ws.Range("Q180").Value() = 25
$ws.Range("R180").Value() = "Hortaliza"
